$d = $word.ActiveDocument

# 1. Add "Pedro Carlos" to the authors list
$d.Content.Find.Execute(
    "Guilherme Bolfe, Silvio Bolfe, Gustavo Bolfe",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Guilherme Bolfe, Silvio Bolfe, Gustavo Bolfe, Pedro Carlos",
    2
)

# 2. Add an extra affiliation placeholder comma in the superscript run
$d.Content.Find.Execute(
    ", , ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", , , ",
    2
)

# 3. Add the new author's e-mail address
$d.Content.Find.Execute(
    "guilhermebolfe11@gmail.com, silviobolfe19@gmail.com, bolfeguilherme@gmail.com",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "guilhermebolfe11@gmail.com, silviobolfe19@gmail.com, bolfeguilherme@gmail.com, pedro@gmail.com",
    2
)

# 4. Flip the exclusion rule from "Any" to "All"
$d.Content.Find.Execute(
    "Exclusion Rule: Any.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Exclusion Rule: All.",
    2
)
